$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$tr.Text = "ENTERPRISE .NET I  "
$tr.InsertAfter(" ")
